$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the group-work (gemeenschappelijk) time for week 1 (column D, row 10)
# from 9:15 to 11:45 (added 2.5 hours of group work on Saturday).
$ws.Range("D10").Value = 11.75 / 24

# Recalculate dependent formulas (D12:D15 etc. depend on D10).
$excel.Calculate()

# Move/restore the active cell selection to D11, as recorded after the edit.
$ws.Range("D11").Select()
